# Rename the existing "strategy_id-5008" sheet to "strategy_id-5007", then
# duplicate it (preserving all data/formatting) and name the copy
# "strategy_id-5009", inserting it right after "strategy_id-5007".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("strategy_id-5008")
$ws.Name = "strategy_id-5007"

# Copy the sheet so the new sheet is placed immediately after it, then
# rename the newly created copy.
$ws.Copy($null, $ws)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "strategy_id-5009"

# Keep the originally active sheet selected (first sheet), since copying
# a worksheet activates the new copy as a side effect.
$wb.Worksheets.Item(1).Activate()
